$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 2: Units -> "g C per individual", Uncertainty -> 2.33830323305659
$ws.Range("C2").Value = "g C per individual"
$ws.Range("D2").Value = 2.33830323305659

# Update row 3: Units -> "Number of individuals", Uncertainty -> 13.610421108098
$ws.Range("C3").Value = "Number of individuals"
$ws.Range("D3").Value = 13.610421108098

# Remove row 4 entirely (the "Total number of non-deep subsurface phages" row)
$ws.Rows.Item(4).Delete()
